# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# 1) Filas 16-27 (ARELIS ESTHER CASTRO ORTIZ): la columna "Periodo Mora" (E)
#    queda en orden ascendente 2001..2012 (antes estaba en orden
#    descendente 2012..2001), y "Salario Basico" (G) sube de 828116 a 877803.
# 2) Filas 28-36 (GINA PAOLA GARCIA MARTINEZ): la columna "Periodo Mora" (E)
#    queda en orden ascendente 2203..2211 (antes estaba en orden
#    descendente 2211..2203), y se intercambian los valores de
#    "Valor Mora" (F) entre la primera fila (28) y la ultima (36).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Bloque 1: filas 16-27 ---------------------------------------------
$startRow1 = 16
$endRow1   = 27

for ($i = 0; $i -lt ($endRow1 - $startRow1 + 1); $i++) {
    $row = $startRow1 + $i
    $year = 2001 + $i
    $ws.Cells.Item($row, 5).Value = [string]$year
    $ws.Cells.Item($row, 7).Value = 877803
}

# --- Bloque 2: filas 28-36 ---------------------------------------------
$startRow2 = 28
$endRow2   = 36

for ($i = 0; $i -lt ($endRow2 - $startRow2 + 1); $i++) {
    $row = $startRow2 + $i
    $year = 2203 + $i
    $ws.Cells.Item($row, 5).Value = [string]$year
}

# Intercambia Valor Mora entre la fila 28 y la fila 36 (28000 <-> 36000)
$ws.Cells.Item(28, 6).Value = 36000
$ws.Cells.Item(36, 6).Value = 28000
